$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "QWE"

$ws.Range("AB2").Value = "Tapped"
$ws.Range("AC2").Value = "Manager seems to be friendly"

$ws.Range("AB3").Value = "LetterGiven"
$ws.Range("AC3").Value = "Manager is not friendly"

$ws.Columns.Item(28).ColumnWidth = 13.8
$ws.Columns.Item(29).ColumnWidth = 8

$ws.Range("E3").Select() | Out-Null

$excel.ActiveWindow.Zoom = 120
